$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1624.4286
$ws.Range("I41").Value = 1535.6666
$ws.Range("J41").Value = 1742.7778
$ws.Range("K41").Value = 1535.6666
$ws.Range("L41").Value = 1742.7778
$ws.Range("M41").Value = -1095.6666
$ws.Range("N41").Value = -2622.7778

$ws.Range("H92").Value = 1045.1578
$ws.Range("I92").Value = 1045.1578
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1045.1578
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 202.8422
$ws.Range("N92").ClearContents()

$ws.Range("H97").Value = 7485.6665
$ws.Range("J97").Value = 7485.6665
$ws.Range("L97").Value = 22456.9995
$ws.Range("N97").Value = -23448.9995

$ws.Range("H116").Value = 8581.4
$ws.Range("I116").Value = 9118.166999999999
$ws.Range("J116").Value = 7776.25
$ws.Range("K116").Value = 9118.166999999999
$ws.Range("L116").Value = 7776.25
$ws.Range("M116").Value = -5676.166999999999
$ws.Range("N116").Value = -14660.25

$ws.Range("H132").Value = 3288.8254
$ws.Range("I132").Value = 3187.4182
$ws.Range("K132").Value = 9562.2546
$ws.Range("M132").Value = -7032.2546

$ws.Range("H138").Value = 4734.5835
$ws.Range("I138").Value = 4156.7
$ws.Range("K138").Value = 12470.1
$ws.Range("M138").Value = -7330.099999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1515
$ws.Range("I2").Value = 899.5
$ws.Range("J2").Value = 1822.75
$ws.Range("K2").Value = 899.5
$ws.Range("L2").Value = 1822.75
$ws.Range("M2").Value = -786.5
$ws.Range("N2").Value = -2048.75

$ws.Range("H32").Value = 1638.6774
$ws.Range("I32").Value = 1327.5747
$ws.Range("K32").Value = 1327.5747
$ws.Range("M32").Value = -1040.5747

$ws.Range("H45").Value = 1684.5
$ws.Range("I45").Value = 1487.2307
$ws.Range("J45").Value = 2197.4
$ws.Range("K45").Value = 1487.2307
$ws.Range("L45").Value = 2197.4
$ws.Range("M45").Value = -1110.2307
$ws.Range("N45").Value = -2951.4

$ws.Range("H116").Value = 1515
$ws.Range("I116").Value = 899.5
$ws.Range("J116").Value = 1822.75
$ws.Range("K116").Value = 899.5
$ws.Range("L116").Value = 1822.75
$ws.Range("M116").Value = 1394.5
$ws.Range("N116").Value = -6410.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1515
$ws.Range("I3").Value = 899.5
$ws.Range("J3").Value = 1822.75
$ws.Range("K3").Value = 899.5
$ws.Range("L3").Value = 1822.75
$ws.Range("M3").Value = -785.5
$ws.Range("N3").Value = -2050.75

$ws.Range("H37").Value = 3500
$ws.Range("I37").Value = 3500
$ws.Range("K37").Value = 3500
$ws.Range("M37").Value = -3363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2493.4546
$ws.Range("I31").Value = 1867.5676
$ws.Range("J31").Value = 3780
$ws.Range("K31").Value = 1867.5676
$ws.Range("L31").Value = 3780
$ws.Range("M31").Value = -1572.5676
$ws.Range("N31").Value = -4370

$ws.Range("H34").Value = 2493.4546
$ws.Range("I34").Value = 1867.5676
$ws.Range("J34").Value = 3780
$ws.Range("K34").Value = 1867.5676
$ws.Range("L34").Value = 3780
$ws.Range("M34").Value = -1665.5676
$ws.Range("N34").Value = -4184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1717.1578
$ws.Range("I26").Value = 165
$ws.Range("J26").Value = 2008.1875
$ws.Range("K26").Value = 495
$ws.Range("L26").Value = 6024.5625
$ws.Range("M26").Value = -207
$ws.Range("N26").Value = -6600.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2677.5
$ws.Range("I80").Value = 2564.1667
$ws.Range("J80").Value = 2762.5
$ws.Range("K80").Value = 2564.1667
$ws.Range("L80").Value = 2762.5
$ws.Range("M80").Value = -1566.1667
$ws.Range("N80").Value = -4758.5

$ws.Range("H83").Value = 2677.5
$ws.Range("I83").Value = 2564.1667
$ws.Range("J83").Value = 2762.5
$ws.Range("K83").Value = 12820.8335
$ws.Range("L83").Value = 13812.5
$ws.Range("M83").Value = -7828.833500000001
$ws.Range("N83").Value = -23796.5

$ws.Range("H113").Value = 5181.125
$ws.Range("I113").Value = 5599.857
$ws.Range("K113").Value = 5599.857
$ws.Range("M113").Value = -3429.857

$ws.Range("H126").Value = 7115
$ws.Range("I126").Value = 4093.8572
$ws.Range("J126").Value = 14164.333
$ws.Range("K126").Value = 12281.5716
$ws.Range("L126").Value = 42492.999
$ws.Range("M126").Value = -9811.571599999999
$ws.Range("N126").Value = -47432.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3201.6667
$ws.Range("I7").Value = 3292.1
$ws.Range("J7").Value = 2749.5
$ws.Range("K7").Value = 3292.1
$ws.Range("L7").Value = 2749.5
$ws.Range("M7").Value = -3180.1
$ws.Range("N7").Value = -2973.5

$ws.Range("H22").Value = 10002
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 10002
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H46").Value = 1687.25
$ws.Range("I46").Value = 1083.8334
$ws.Range("K46").Value = 1083.8334
$ws.Range("M46").Value = -895.8334

$ws.Range("H61").Value = 1265.6
$ws.Range("I61").Value = 1284.5714
$ws.Range("K61").Value = 1284.5714
$ws.Range("M61").Value = -1082.5714

$ws.Range("H82").Value = 1249
$ws.Range("I82").Value = 1133.1
$ws.Range("K82").Value = 1133.1
$ws.Range("M82").Value = -772.0999999999999

$ws.Range("H85").Value = 1249
$ws.Range("I85").Value = 1133.1
$ws.Range("K85").Value = 1133.1
$ws.Range("M85").Value = 114.9000000000001

$ws.Range("H113").Value = 1265.6
$ws.Range("I113").Value = 1284.5714
$ws.Range("K113").Value = 1284.5714
$ws.Range("M113").Value = 885.4286

$ws.Range("H126").Value = 3201.6667
$ws.Range("I126").Value = 3292.1
$ws.Range("J126").Value = 2749.5
$ws.Range("K126").Value = 9876.299999999999
$ws.Range("L126").Value = 8248.5
$ws.Range("M126").Value = -7406.299999999999
$ws.Range("N126").Value = -13188.5

$ws.Range("H136").Value = 1505.7632
$ws.Range("I136").Value = 1019.5926
$ws.Range("K136").Value = 3058.7778
$ws.Range("M136").Value = -508.7777999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 778.2
$ws.Range("I107").Value = 672.75
$ws.Range("K107").Value = 2018.25
$ws.Range("M107").Value = -98.25
